# Fix missed word in the "Changes & Misc." slide bullet:
# "... select with approver ..." -> "... select which approver ..."
#
# The target paragraph is split across 3 runs whose text boundaries shift
# (run 1 absorbs the old run-2 text, the "select which " phrase becomes its
# own run, and the trailing "approver they would like." becomes the third
# run) even though only a single word actually changed. Setting each run's
# TextRange.Text individually (rather than the paragraph/shape text as a
# whole) performs a direct, non-diffed replacement of that run's content,
# which lets us reproduce the exact same run layout as the authored edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(3)

$r1 = $para.Runs(1)
$r2 = $para.Runs(2)
$r3 = $para.Runs(3)

$r1.Text = "When requesting approval, the request will either automatically be sent to any of the approvers for that software or the user can manually "
$r2.Text = "select which "
$r3.Text = "approver they would like."
